{"js": "// Office.js (Word JavaScript API) script.\n// Applies the textual edits described by the diff:\n//  - Updates the page title/heading text (two occurrences: the H1 heading\n//    and the bold \"title\" run near the end of the document).\n//  - Expands several \"What we like\" / \"What we don't like\" bullet points.\n//  - Updates the italic meta-description run near the end of the document.\n\n// NOTE: order matters. The meta-description sentence (originally ending in\n// \"...Super Wild system enhances winning chances.\") fully contains the\n// shorter bullet-point phrase \"Super Wild system enhances winning chances\".\n// Replacing the long sentence first avoids the shorter search term\n// incorrectly matching (and only partially updating) the long sentence.\nconst replacements = [\n  {\n    search: \"Detailed review of Fire Archer slot machine. Play for free with beautiful graphics and innovative gameplay. Super Wild system enhances winning chances.\",\n    replace: \"Detailed review of Fire Archer slot game with stunning graphics and Super Wild system. Play for free and win big.\",\n  },\n  {\n    search: \"Play Fire Archer Free: Stunning Graphics & Super Wild System\",\n    replace: \"Play Fire Archer Free - Review of Fire Archer Slot Game\",\n  },\n  {\n    search: \"Stunning graphics\",\n    replace: \"Stunning graphics and animation-like visuals\",\n  },\n  {\n    search: \"Super Wild system enhances winning chances\",\n    replace: \"Super Wild system enhances chances of winning\",\n  },\n  {\n    search: \"Two special symbols and features\",\n    replace: \"Possibility to play for free in demo mode\",\n  },\n  {\n    search: \"Good theoretical RTP\",\n    replace: \"Good theoretical Return to Player (RTP) of 96.07%\",\n  },\n  {\n    search: \"High volatility can make it challenging to win\",\n    replace: \"High volatility may lead to infrequent wins\",\n  },\n  {\n    search: \"Substantial payout potential may not appeal to all players\",\n    replace: \"Max potential payout of 4,000 times the stake can be challenging to obtain\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { search, replace } of replacements) {\n  const results = body.search(search, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM (PowerShell-style) script.\n# Applies the textual edits described by the diff:\n#  - Updates the page title/heading text (two occurrences: the H1 heading\n#    and the bold \"title\" run near the end of the document).\n#  - Expands several \"What we like\" / \"What we don't like\" bullet points.\n#  - Updates the italic meta-description run near the end of the document.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# NOTE: order matters. The meta-description sentence (originally ending in\n# \"...Super Wild system enhances winning chances.\") fully contains the\n# shorter bullet-point phrase \"Super Wild system enhances winning chances\".\n# Replacing the long sentence first avoids the shorter search term\n# incorrectly matching (and only partially updating) the long sentence.\nReplace-AllText \"Detailed review of Fire Archer slot machine. Play for free with beautiful graphics and innovative gameplay. Super Wild system enhances winning chances.\" \"Detailed review of Fire Archer slot game with stunning graphics and Super Wild system. Play for free and win big.\"\nReplace-AllText \"Play Fire Archer Free: Stunning Graphics & Super Wild System\" \"Play Fire Archer Free - Review of Fire Archer Slot Game\"\nReplace-AllText \"Stunning graphics\" \"Stunning graphics and animation-like visuals\"\nReplace-AllText \"Super Wild system enhances winning chances\" \"Super Wild system enhances chances of winning\"\nReplace-AllText \"Two special symbols and features\" \"Possibility to play for free in demo mode\"\nReplace-AllText \"Good theoretical RTP\" \"Good theoretical Return to Player (RTP) of 96.07%\"\nReplace-AllText \"High volatility can make it challenging to win\" \"High volatility may lead to infrequent wins\"\nReplace-AllText \"Substantial payout potential may not appeal to all players\" \"Max potential payout of 4,000 times the stake can be challenging to obtain\"\n"}
